# -----------------------------------------------------------------------
# Add a new "2022-Q4" worksheet (right after "总计") with its fund detail
# data, and update the "总计" (summary) sheet with a new leading row for
# 2022-Q4 totals, shifting all the other rows down by one.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# -------------------------------------------------------------------
# 1. Update "总计" summary sheet: insert a new 2022-Q4 row at the top
#    of the data (row 2), pushing the existing rows down by one.
# -------------------------------------------------------------------

# Shift rows 9..2 down to 10..3 (process bottom-up so we never clobber
# a row before we've read it).
for ($r = 9; $r -ge 2; $r--) {
    $dst = $r + 1
    $ws1.Cells.Item($dst, 1).Value = $r - 1
    $ws1.Cells.Item($dst, 2).Value2 = $ws1.Cells.Item($r, 2).Value2
    $ws1.Cells.Item($dst, 3).Value2 = $ws1.Cells.Item($r, 3).Value2
    $ws1.Cells.Item($dst, 4).Value2 = $ws1.Cells.Item($r, 4).Value2
}

# The freshly created A10 cell needs the same box/bold formatting that
# the rest of column A uses; copy it from the row above.
$ws1.Cells.Item(9, 1).Copy()
$ws1.Cells.Item(10, 1).PasteSpecial(-4122)

# Fill in the new 2022-Q4 summary row.
$ws1.Cells.Item(2, 1).Value = 0
$ws1.Cells.Item(2, 2).Value = "2022-Q4"
$ws1.Cells.Item(2, 3).Value = 11
$ws1.Cells.Item(2, 4).Value = 5.5

# -------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet right after "总计" and fill it
#    with the quarterly fund holding detail.
# -------------------------------------------------------------------

$wsNew = $wb.Worksheets.Add($null, $ws1)
$wsNew.Name = "2022-Q4"

# Borrow formatting (bold + border + centered) from an existing detail
# sheet so the new sheet matches the look of its siblings.
$wsStyleSrc = $wb.Worksheets.Item("2022-Q3")

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $wsNew.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$rows = @(
    @(0, "003835", "鹏华沪深港新兴成长灵活配置混合", "45.58", "94.28", "3.83", "1.7457", 6),
    @(1, "004854", "广发中证全指汽车指数A", "15.08", "92.79", "5.13", "0.7736", 7),
    @(2, "016067", "鹏华新能源汽车混合A", "19.04", "95.25", "3.89", "0.7407", 5),
    @(3, "010094", "交银施罗德产业机遇混合", "16.65", "85.76", "3.79", "0.6310", 7),
    @(4, "004855", "广发中证全指汽车指数C", "11.13", "92.79", "5.13", "0.5710", 7),
    @(5, "519773", "交银施罗德数据产业灵活配置混合A", "13.79", "86.66", "3.95", "0.5447", 7),
    @(6, "398061", "中海消费混合", "3.80", "90.62", "4.67", "0.1775", 2),
    @(7, "016068", "鹏华新能源汽车混合C", "3.65", "95.25", "3.89", "0.1420", 5),
    @(8, "014549", "交银施罗德数据产业灵活配置混合C", "2.43", "86.66", "3.95", "0.0960", 7),
    @(9, "015986", "中海新兴成长六个月持有期混合", "2.56", "40.87", "2.86", "0.0732", 2),
    @(10, "510770", "申万菱信上证G60战略新兴产业成份ETF", "0.23", "95.96", "3.32", "0.0076", 10)
)

$r = 2
foreach ($row in $rows) {
    $wsNew.Cells.Item($r, 1).Value = $row[0]
    $wsNew.Cells.Item($r, 2).Value = "'" + $row[1]
    $wsNew.Cells.Item($r, 3).Value = $row[2]
    $wsNew.Cells.Item($r, 4).Value = "'" + $row[3]
    $wsNew.Cells.Item($r, 5).Value = "'" + $row[4]
    $wsNew.Cells.Item($r, 6).Value = "'" + $row[5]
    $wsNew.Cells.Item($r, 7).Value = "'" + $row[6]
    $wsNew.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Apply header-row formatting (bold font + border + centered alignment).
$wsStyleSrc.Range("B1:H1").Copy()
$wsNew.Range("B1:H1").PasteSpecial(-4122)

# Apply column-A formatting for every data row (same boxed/bold style
# used by column A elsewhere in the workbook).
$wsStyleSrc.Range("A2").Copy()
$wsNew.Range("A2:A12").PasteSpecial(-4122)

# -------------------------------------------------------------------
# 3. Restore the originally-active tab (2020-Q4), since adding the new
#    sheet above switched the active tab to it.
# -------------------------------------------------------------------
$wsActive = $wb.Worksheets.Item("2020-Q4")
$wsActive.Activate()

Write-Output "done"
